$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "punkt_adt": update ydt_light / ydt_heavy for row 31
# ---------------------------------------------------------------------------
$wsAdt = $wb.Worksheets.Item("punkt_adt")
$wsAdt.Range("P31").Value = 10070
$wsAdt.Range("Q31").Value = 1420

# ---------------------------------------------------------------------------
# Sheet "punktindeks_maned": add "nov" (column O) figures for the TRPs that
# now have a November reading
# ---------------------------------------------------------------------------
$wsMnd = $wb.Worksheets.Item("punktindeks_maned")
$wsMnd.Range("O13").Value = -1
$wsMnd.Range("O16").Value = 0.47
$wsMnd.Range("O19").Value = 4.7
$wsMnd.Range("O25").Value = -8.73
$wsMnd.Range("O34").Value = 2.17
$wsMnd.Range("O40").Value = -0.24
$wsMnd.Range("O49").Value = -0.45
$wsMnd.Range("O58").Value = -7.43
$wsMnd.Range("O61").Value = 6.39
$wsMnd.Range("O64").Value = 0.54
$wsMnd.Range("O73").Value = -5.31
$wsMnd.Range("O76").Value = 9.630000000000001
$wsMnd.Range("O79").Value = 8.83

# ---------------------------------------------------------------------------
# Sheet "byindeks_aarlig": extend the direct 2024-2025 index by one month
# (jan-okt -> jan-nov)
# ---------------------------------------------------------------------------
$wsAar = $wb.Worksheets.Item("byindeks_aarlig")
$wsAar.Range("C5").Value = 11
$wsAar.Range("F5").Value = "jan-nov"
$wsAar.Range("H5").Value = 1.0127
$wsAar.Range("I5").Value = 1.3
$wsAar.Range("J5").Value = 1.141741421857783
$wsAar.Range("K5").Value = -0.9
$wsAar.Range("L5").Value = 3.5

# ---------------------------------------------------------------------------
# Sheet "by_glid_indeks": a new rolling 12-month window (des 2024 - nov 2025)
# is inserted right before the existing 24-month series (which all shift
# down by one row), and a new rolling 24-month window (des 2023 - nov 2025)
# is appended at the end.
# ---------------------------------------------------------------------------
$wsGlid = $wb.Worksheets.Item("by_glid_indeks")

# Insert a fresh row at 25, pushing the rest of the 24_months block down.
$wsGlid.Rows.Item(25).Insert()

$wsGlid.Range("A25").Value = 0.956676435683516
$wsGlid.Range("B25").Value = -4.332356431648399
$wsGlid.Range("C25").Value = 14
$wsGlid.Range("D25").Value = 9.461194469994499
$wsGlid.Range("E25").Value = 10.84107229899311
$wsGlid.Range("F25").Value = 3.487192733000519
$wsGlid.Range("G25").Value = -11.2
$wsGlid.Range("H25").Value = 2.5
$wsGlid.Range("I25").Value = "2019 - (des 2024 - nov 2025)"
$wsGlid.Range("J25").Value = 45962
$wsGlid.Range("J25").NumberFormat = "yyyy-mm-dd"
$wsGlid.Range("K25").Value = 11
$wsGlid.Range("L25").Value = 2025
$wsGlid.Range("M25").Value = "12_months"

# Append the new trailing 24-month row (row 37, after the insert above
# shifted the previous last row 35 down to 36).
$wsGlid.Range("A37").Value = 0.9505415319810241
$wsGlid.Range("B37").Value = -4.945846801897591
$wsGlid.Range("C37").Value = 14
$wsGlid.Range("D37").Value = 9.461194469994499
$wsGlid.Range("E37").Value = 9.649200500987547
$wsGlid.Range("F37").Value = 3.1275596420733
$wsGlid.Range("G37").Value = -11.1
$wsGlid.Range("H37").Value = 1.2
$wsGlid.Range("I37").Value = "2019 - (des 2023 - nov 2025)"
$wsGlid.Range("J37").Value = 45962
$wsGlid.Range("J37").NumberFormat = "yyyy-mm-dd"
$wsGlid.Range("K37").Value = 11
$wsGlid.Range("L37").Value = 2025
$wsGlid.Range("M37").Value = "24_months"
